# PlayerPerformance_5985.xlsx update:
#   1. Add a new worksheet "ODI Batting Extra" (sheetId 4) at the end of the workbook.
#   2. Populate it with per-match batting-extra stats keyed by MATCH_CODE.
#   3. On the existing "ODI Batting" sheet, clear out the (previously blank)
#      INNING_NUMBER cells for the three "did not bat" rows (B2, B3, B7).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the new sheet after the last existing one so it lands at the end
# ---------------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "ODI Batting Extra"

# Helper-less: set a cell as real text (avoids Excel's auto "looks like a
# number/percent" coercion) and strip the quote-prefix style it leaves behind
# so only genuine header cells keep a non-default style.
function Set-Text($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

function Set-BlankText($range) {
    $range.Value = "'"
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------------
# 2. Header row (bold / bordered / centered, same look as the other sheets)
# ---------------------------------------------------------------------------
$ws.Range("A1").Value = "MATCH_CODE"
$ws.Range("B1").Value = "BATTING_POSITION"
$ws.Range("C1").Value = "NUM_4"
$ws.Range("D1").Value = "NUM_6"
$ws.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$ws.Range("F1").Value = "MAN_OF_MATCH"

# Reuse the exact header style already used on every other sheet (bold,
# thin border all round, centered/top aligned) instead of fabricating a new
# (near-duplicate) style entry.
$headerSample = $wb.Worksheets.Item("ODI Batting").Range("A1")
$headerSample.Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 3. Data rows
# ---------------------------------------------------------------------------
Set-Text $ws.Range("A2") "4517"
Set-BlankText $ws.Range("B2")
Set-BlankText $ws.Range("C2")
Set-BlankText $ws.Range("D2")
Set-BlankText $ws.Range("E2")
Set-Text $ws.Range("F2") "NO"

Set-Text $ws.Range("A3") "4526"
$ws.Range("B3").Value = 9
Set-BlankText $ws.Range("C3")
Set-BlankText $ws.Range("D3")
Set-BlankText $ws.Range("E3")
Set-Text $ws.Range("F3") "NO"

Set-Text $ws.Range("A4") "4529"
$ws.Range("B4").Value = 10
Set-Text $ws.Range("C4") "0"
Set-Text $ws.Range("D4") "0"
Set-BlankText $ws.Range("E4")
Set-Text $ws.Range("F4") "NO"

Set-Text $ws.Range("A5") "4698"
$ws.Range("B5").Value = 8
Set-Text $ws.Range("C5") "1"
Set-Text $ws.Range("D5") "0"
Set-Text $ws.Range("E5") "1.68%"
Set-Text $ws.Range("F5") "YES"

Set-Text $ws.Range("A6") "4700"
$ws.Range("B6").Value = 9
Set-Text $ws.Range("C6") "0"
Set-Text $ws.Range("D6") "0"
Set-Text $ws.Range("E6") "0.70%"
Set-Text $ws.Range("F6") "NO"

Set-Text $ws.Range("A7") "4746"
$ws.Range("B7").Value = 8
Set-BlankText $ws.Range("C7")
Set-BlankText $ws.Range("D7")
Set-BlankText $ws.Range("E7")
Set-Text $ws.Range("F7") "YES"

Set-Text $ws.Range("A8") "4751"
Set-BlankText $ws.Range("B8")
Set-BlankText $ws.Range("C8")
Set-BlankText $ws.Range("D8")
Set-BlankText $ws.Range("E8")
Set-Text $ws.Range("F8") "NO"

$ws.Range("A1").Select()

# ---------------------------------------------------------------------------
# 4. "ODI Batting": the blank INNING_NUMBER cells for non-batting innings are
#    fully cleared (no cell at all), matching rows 2, 3 and 7.
# ---------------------------------------------------------------------------
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("B2").ClearContents()
$batting.Range("B3").ClearContents()
$batting.Range("B7").ClearContents()
